# Econ 470 grade-calculator update:
# Scale down the raw point values entered for each grading component
# (the underlying assignment/exam points were re-entered with one
# fewer digit), letting the existing formulas recompute weighted
# points and totals automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 15    # LabQuizzes
$ws.Range("B4").Value = 6     # MidtermExam1
$ws.Range("B5").Value = 23    # MidtermExam2
$ws.Range("B6").Value = 23    # HomeworkAssignments
$ws.Range("B7").Value = 33    # Final 330

# Leave the cursor where the author last left it.
$ws.Range("C15").Select() | Out-Null
